$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Day 4 date changed from 29/05/2019 to 26/06/2019, and a new note added in B4
$ws.Range("A4").Value = "Day 4 (26/06/2019)"
$ws.Range("B4").Value = "Learning Converting Our CSS Code to Sass Variables and Nesting => Implement to the Natours Project"

# Widen column B to fit new text (88.5 is the closest achievable input that
# rounds to the target stored width of 89.28515625 given this runtime's
# column-width quantization)
$ws.Columns.Item(2).ColumnWidth = 88.5

# Update selection to B6 (last cell touched)
$ws.Range("B6").Select()
